# Corrections mineures dans le tableau "cas-2".
$d = $word.ActiveDocument

# Etape 2 - Reponses attendues : accorder "affiché" -> "affichées" et
# "le champs" -> "les champs" (accord grammatical, pluriel).
$d.Content.Find.Execute("en question sont affiché dans le champs de texte", $true, $false, $false, $false, $false, $true, 1, $false, "en question sont affichées dans les champs de texte", 2) | Out-Null

# Etape 3 - Reponses attendues : "afficher" -> "affichée" et
# "son champs" -> "son champ" (accord grammatical, singulier).
$d.Content.Find.Execute("La facture est afficher dans son champs de texte", $true, $false, $false, $false, $false, $true, 1, $false, "La facture est affichée dans son champ de texte", 2) | Out-Null
